$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.234005928039551
$ws.Range("B1").Value = 4.983691692352295
$ws.Range("C1").Value = 4.112290859222412
$ws.Range("D1").Value = 4.908996105194092
$ws.Range("E1").Value = 4.686135292053223
